$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 130 (shifts existing rows 130-162 down to 131-163)
$ws.Rows.Item(130).Insert()

# Populate the newly inserted row 130 with a new weekly price record
# (same market/category/quality as the surrounding rows, new date/price)
$ws.Cells.Item(130, 1).Value = 5
$ws.Cells.Item(130, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(130, 3).Value = "Maule"
$ws.Cells.Item(130, 4).Value = [DateTime]"2022-06-10"
$ws.Cells.Item(130, 5).Value = 7
$ws.Cells.Item(130, 6).Value = 100112031
$ws.Cells.Item(130, 7).Value = "Poroto verde"
$ws.Cells.Item(130, 8).Value = "Sin especificar"
$ws.Cells.Item(130, 9).Value = "Primera"
$ws.Cells.Item(130, 10).Value = 100
$ws.Cells.Item(130, 11).Value = 25000
$ws.Cells.Item(130, 12).Value = 25000
$ws.Cells.Item(130, 13).Value = 25000
$ws.Cells.Item(130, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(130, 15).Value = "Región del Maule"
$ws.Cells.Item(130, 16).Value = 1000
$ws.Cells.Item(130, 17).Value = 25
$ws.Cells.Item(130, 18).Value = "Hortaliza"
